$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 28, shifting the existing rows (28-50) down to (29-51).
$ws.Rows(28).Insert()

# Populate the newly inserted row 28 with the "Pre Profile Presentation" / "irrelevant"
# phase row that explains the upcoming matrix-game instructions (male/female variants).
$ws.Range("A28").Value = "Pre Profile Presentation"
$ws.Range("B28").Value = "irrelevant"

# Set E28 (female wording) first, then C28/D28 (male wording) so new shared strings are
# appended in that order.
$ws.Range("E28").Value = "לאחר שתצפי בתכונות והמאפיינים של המשתתף האחר, תוצג בפניך טבלה דומה לזו שהוצגה קודם אשר מתארת משחק בו קיימות שתי אפשרויות בחירה שלך ושתי אפשרויות בחירה של האחר או האחרת. תצטרכי להחליט כיצד לפעול מול האדם אשר בתכונות ובמאפיינים שלו צפית, והאדם האחר יצטרף להחליט כיצד לפעול מולך. אם תרצי להתרשם פעם נוספת לפני החלטתך מהתכונות והמאפיינים של האדם האחר, תוכלי ללחוץ על הכפתור ""הצג שוב"". אחרי שתבחרי באחת האפשרויות, יוצג בפניך אדם אחר."
$ws.Range("C28").Value = "לאחר שתצפה בתכונות והמאפיינים של המשתתף האחר, תוצג בפניך טבלה דומה לזו שהוצגה קודם אשר מתארת משחק בו קיימות שתי אפשרויות בחירה שלך ושתי אפשרויות בחירה של האחר או האחרת. תצטרך להחליט כיצד לפעול מול האדם אשר בתכונות ובמאפיינים שלו צפית, והאדם האחר יצטרף להחליט כיצד לפעול מולך. אם תרצה להתרשם פעם נוספת לפני החלטתך מהתכונות והמאפיינים של האדם האחר, תוכל ללחוץ על הכפתור ""הצג שוב"". אחרי שתבחר באחת האפשרויות, יוצג בפניך אדם אחר."
$ws.Range("D28").Value = "לאחר שתצפה בתכונות והמאפיינים של המשתתף האחר, תוצג בפניך טבלה דומה לזו שהוצגה קודם אשר מתארת משחק בו קיימות שתי אפשרויות בחירה שלך ושתי אפשרויות בחירה של האחר או האחרת. תצטרך להחליט כיצד לפעול מול האדם אשר בתכונות ובמאפיינים שלו צפית, והאדם האחר יצטרף להחליט כיצד לפעול מולך. אם תרצה להתרשם פעם נוספת לפני החלטתך מהתכונות והמאפיינים של האדם האחר, תוכל ללחוץ על הכפתור ""הצג שוב"". אחרי שתבחר באחת האפשרויות, יוצג בפניך אדם אחר."

# Match the wrapped-text row height used for this instructions row.
$ws.Rows(28).RowHeight = 129.6

# Update the saved view so it scrolls/selects the newly inserted row.
$ws.Range("A28").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
